$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.7332147692668679
$ws.Range("D2").Value = 0.4684530523185102

# Row 3
$ws.Range("C3").Value = 1.865406320549269
$ws.Range("D3").Value = 0.07077305411015056
$ws.Range("G3").Value = "No"

# Row 4
$ws.Range("C4").Value = 0.5901358428226466
$ws.Range("D4").Value = 0.5590016776715054

# Row 5
$ws.Range("C5").Value = 2.126580489148894
$ws.Range("D5").Value = 0.04079462919154664
$ws.Range("G5").Value = "Sí"

# Row 6
$ws.Range("C6").Value = 1.845686814536893
$ws.Range("D6").Value = 0.07366069461930147

# Row 7
$ws.Range("C7").Value = -0.03727037138849371
$ws.Range("D7").Value = 0.9704874042606013

# Row 8
$ws.Range("C8").Value = 1.93053244498202
$ws.Range("D8").Value = 0.06191546590421226

# Row 9
$ws.Range("C9").Value = -1.220976153158016
$ws.Range("D9").Value = 0.2304948034534746

# Row 10
$ws.Range("C10").Value = -0.1934813693309636
$ws.Range("D10").Value = 0.8477332530457553

# Row 11
$ws.Range("C11").Value = 1.402335401157454
$ws.Range("D11").Value = 0.169880032411593
